# Update the "ltp" sheet with refreshed LTP/PREV quote data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")
$ws.Activate()

$ws.Range("B2").Value = 675.75
$ws.Range("C2").Value = 675.25
$ws.Range("B3").Value = 3261.75
$ws.Range("C3").Value = 3224.8
$ws.Range("B4").Value = 476.45
$ws.Range("C4").Value = 476.7
$ws.Range("B5").Value = 1591.75
$ws.Range("C5").Value = 1590.1
$ws.Range("B6").Value = 7341.55
$ws.Range("C6").Value = 7294.95
$ws.Range("B7").Value = 194.34
$ws.Range("C7").Value = 194.23
$ws.Range("B8").Value = 286.25
$ws.Range("C8").Value = 282.7
$ws.Range("B9").Value = 50042.4
$ws.Range("C9").Value = 49912.65
$ws.Range("B10").Value = 878.6
$ws.Range("C10").Value = 874.05
$ws.Range("B11").Value = 4935.1
$ws.Range("C11").Value = 4801.5
$ws.Range("B12").Value = 174.4
$ws.Range("C12").Value = 172.26
$ws.Range("B13").Value = 1431.05
$ws.Range("C13").Value = 1444.15
$ws.Range("B14").Value = 683.6
$ws.Range("C14").Value = 680.7
$ws.Range("B15").Value = 1502.35
$ws.Range("C15").Value = 1507.25
$ws.Range("B16").Value = 1052.45
$ws.Range("C16").Value = 1041.25
$ws.Range("B17").Value = 731.65
$ws.Range("C17").Value = 719.3
$ws.Range("B18").Value = 2928.6
$ws.Range("C18").Value = 2861.7
$ws.Range("B19").Value = 298.95
$ws.Range("C19").Value = 293.3
$ws.Range("B20").Value = 23466.65
$ws.Range("C20").Value = 23399.25
$ws.Range("B21").Value = 368.45
$ws.Range("C21").Value = 369.95
$ws.Range("B22").Value = 839.2
$ws.Range("C22").Value = 843.9
$ws.Range("B23").Value = 757.55
$ws.Range("C23").Value = 759.3
$ws.Range("B24").Value = 993.4
$ws.Range("C24").Value = 985.85
$ws.Range("B25").Value = 448.65
$ws.Range("C25").Value = 451.9
$ws.Range("B26").Value = 183.15
$ws.Range("C26").Value = 182.56

# Move the selection to match the post-edit cursor position.
$ws.Range("I9").Select()
